# "Generate Report for Handback"
#
# The localization-status report gets refreshed for a handback cycle:
#   * Every "Ready for handoff" status cell (Overview + per-language sheets)
#     flips to "Handed back: in sync with en-US".
#   * Each per-language sheet (zh-cn, de-de) grows two new columns of data
#     for rows 2-3: "Latest Target File" (E) and "Latest Handback File" (F),
#     both populated as hyperlinked file names, mirroring the existing
#     "Source File Name" (A) / "Latest Handoff File" (C) hyperlinks.
#   * The "Latest Handback DateTime" column (G) for rows 2-3 is stamped
#     with the handback timestamp for that language.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("B3").Value = $newStatus

$zhMdFile = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$zhXlfFile = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"
$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5ac7683945ed341faaf462ccc9a400a1a8126d35/e2e/48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/41e0d1f490bb3de1c1ff71678f7497aae115e44b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), $zhMdUrl, "", "", $zhMdFile)
$wsZhCn.Range("E2").Font.Underline = 2
$wsZhCn.Range("E2").Font.Color = 0xED9564

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhXlfUrl, "", "", $zhXlfFile)
$wsZhCn.Range("F2").Font.Underline = 2
$wsZhCn.Range("F2").Font.Color = 0xED9564

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), $zhMdUrl, "", "", $zhMdFile)
$wsZhCn.Range("E3").Font.Underline = 2
$wsZhCn.Range("E3").Font.Color = 0xED9564

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhXlfUrl, "", "", $zhXlfFile)
$wsZhCn.Range("F3").Font.Underline = 2
$wsZhCn.Range("F3").Font.Color = 0xED9564

$wsZhCn.Range("G2").Value = "2016-02-22 05:08:50"
$wsZhCn.Range("G3").Value = "2016-02-22 05:08:50"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("B3").Value = $newStatus

$deMdFile = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$deXlfFile = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"
$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5ac7683945ed341faaf462ccc9a400a1a8126d35/e2e/48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2987ccb40f6f9233e2a195117bbe91be24e02410/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), $deMdUrl, "", "", $deMdFile)
$wsDeDe.Range("E2").Font.Underline = 2
$wsDeDe.Range("E2").Font.Color = 0xED9564

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deXlfUrl, "", "", $deXlfFile)
$wsDeDe.Range("F2").Font.Underline = 2
$wsDeDe.Range("F2").Font.Color = 0xED9564

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), $deMdUrl, "", "", $deMdFile)
$wsDeDe.Range("E3").Font.Underline = 2
$wsDeDe.Range("E3").Font.Color = 0xED9564

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deXlfUrl, "", "", $deXlfFile)
$wsDeDe.Range("F3").Font.Underline = 2
$wsDeDe.Range("F3").Font.Color = 0xED9564

$wsDeDe.Range("G2").Value = "2016-02-22 05:09:15"
$wsDeDe.Range("G3").Value = "2016-02-22 05:09:15"
